$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.505.91"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.828.23"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.72"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5176"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3872"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08299"
$ws.Range("E9").Value = "  +7.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.123"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.407"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.505"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "1.827.87"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.25"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001120"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.83"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.069"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "28.550.42"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.260"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.44"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "2.036.55"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.414"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.26"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1095"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.101"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.744"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07600"
$ws.Range("E34").Value = "  +7.70%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2236"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02375"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.301"
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.88"
$ws.Range("E39").Value = "  +6.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.781"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6403"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.64"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6169"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.807"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.65"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.003"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.206"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06986"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.38"
$ws.Range("E51").Value = "  +0.46%  "
